{"js": "// Original paragraph text (right after the italic word \"unscrambling\"):\n//   \"For unscrambling you should follow similar steps as above:\"\n// Target paragraph text:\n//   \"For unscrambling, follow the steps as below:\"\n// The new leading \",\" right after \"unscrambling\" keeps the italic\n// formatting (it becomes part of the same run as \"unscrambling\"); the\n// rest of the new text (\"  follow \", \"the \", \"steps as below\") is plain\n// (non-italic), matching the style of the text it replaces.\n\nconst body = context.document.body;\n\n// 1) Locate the old tail sentence and replace it with the new tail text.\n//    Leaving this as a single insertText(\"Replace\") call keeps the run's\n//    existing (non-italic) formatting, so no explicit italic reset needed.\nconst oldTail = body.search(\" you should follow similar steps as above\", {\n  matchCase: true\n});\noldTail.load(\"items\");\nawait context.sync();\n\nif (oldTail.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to replace.\");\n}\n\noldTail.items[0].insertText(\" follow the steps as below\", \"Replace\");\nawait context.sync();\n\n// 2) Re-locate the freshly inserted text and prepend an italic \",\" right\n//    before it (i.e. immediately after \"unscrambling\").\nconst newTail = body.search(\" follow the steps as below\", {\n  matchCase: true\n});\nnewTail.load(\"items\");\nawait context.sync();\n\nif (newTail.items.length === 0) {\n  throw new Error(\"Could not find the replacement text to prepend the comma to.\");\n}\n\nconst commaRange = newTail.items[0].insertText(\",\", \"Before\");\ncommaRange.font.italic = true;\n\nawait context.sync();\n", "ps1": "# Original paragraph text (right after the italic word \"unscrambling\"):\n#   \"For unscrambling you should follow similar steps as above:\"\n# Target paragraph text:\n#   \"For unscrambling, follow the steps as below:\"\n# The new leading \",\" right after \"unscrambling\" keeps the italic\n# formatting (it joins the same run as \"unscrambling\"); the rest of the\n# new text (\"  follow \", \"the \", \"steps as below\") stays plain\n# (non-italic), matching the style of the text it replaces.\n\n$d = $word.ActiveDocument\n\n# 1) Find the old tail sentence and replace it with the new tail text.\n#    Reusing the matched Range's own .Text setter keeps the run's existing\n#    (non-italic) formatting, so no explicit italic reset is required.\n$oldTail = $d.Content\n$found = $oldTail.Find.Execute(\" you should follow similar steps as above\")\nif (-not $found) {\n    throw \"Could not find the target sentence to replace.\"\n}\n$oldTail.Text = \" follow the steps as below\"\n\n# 2) Re-locate the freshly inserted text and insert an italic \",\" right\n#    before it (i.e. immediately after \"unscrambling\").\n$newTail = $d.Content\n$found2 = $newTail.Find.Execute(\" follow the steps as below\")\nif (-not $found2) {\n    throw \"Could not find the replacement text to prepend the comma to.\"\n}\n$commaRange = $newTail.Duplicate\n$commaRange.Collapse(1)   # wdCollapseStart\n$commaRange.InsertBefore(\",\")\n$commaRange.Font.Italic = $true\n"}
